$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last row (row 53) - the dataset now spans one fewer row
$ws.Rows(53).Delete()

# Recalculated naive forecaster values for rows 2-52 (bugfix: correct date/value alignment)
$ws.Cells.Item(2, "A").Value = 39583
$ws.Cells.Item(2, "B").Value = 2008
$ws.Cells.Item(2, "C").ClearContents() | Out-Null
$ws.Cells.Item(2, "D").Value = 2009
$ws.Cells.Item(2, "E").Value = 0.5087393606160395
$ws.Cells.Item(3, "A").Value = 39765
$ws.Cells.Item(3, "B").Value = 2008
$ws.Cells.Item(3, "C").ClearContents() | Out-Null
$ws.Cells.Item(3, "D").Value = 2009
$ws.Cells.Item(3, "E").Value = -0.4513776153963867
$ws.Cells.Item(4, "A").Value = 39948
$ws.Cells.Item(4, "B").Value = 2009
$ws.Cells.Item(4, "C").Value = -1.118515468742087
$ws.Cells.Item(4, "D").Value = 2010
$ws.Cells.Item(4, "E").Value = -0.6296678961043134
$ws.Cells.Item(5, "A").Value = 40130
$ws.Cells.Item(5, "B").Value = 2009
$ws.Cells.Item(5, "C").Value = -1.324983933426882
$ws.Cells.Item(5, "D").Value = 2010
$ws.Cells.Item(5, "E").Value = -0.8803581938132576
$ws.Cells.Item(6, "A").Value = 40310
$ws.Cells.Item(6, "B").Value = 2010
$ws.Cells.Item(6, "C").Value = -0.1156872058426073
$ws.Cells.Item(6, "D").Value = 2011
$ws.Cells.Item(6, "E").Value = -0.5120992642018263
$ws.Cells.Item(7, "A").Value = 40494
$ws.Cells.Item(7, "B").Value = 2010
$ws.Cells.Item(7, "C").Value = -0.3900454704678369
$ws.Cells.Item(7, "D").Value = 2011
$ws.Cells.Item(7, "E").Value = -1.213027585730386
$ws.Cells.Item(8, "A").Value = 40676
$ws.Cells.Item(8, "B").Value = 2011
$ws.Cells.Item(8, "C").Value = -0.4084169314491404
$ws.Cells.Item(8, "D").Value = 2012
$ws.Cells.Item(8, "E").Value = -0.6403426624573716
$ws.Cells.Item(9, "A").Value = 40862
$ws.Cells.Item(9, "B").Value = 2011
$ws.Cells.Item(9, "C").Value = -0.2995848153489522
$ws.Cells.Item(9, "D").Value = 2012
$ws.Cells.Item(9, "E").Value = -0.3230872999110068
$ws.Cells.Item(10, "A").Value = 41044
$ws.Cells.Item(10, "B").Value = 2012
$ws.Cells.Item(10, "C").Value = -0.2188016966516937
$ws.Cells.Item(10, "D").Value = 2013
$ws.Cells.Item(10, "E").Value = -0.1561757764150462
$ws.Cells.Item(11, "A").Value = 41228
$ws.Cells.Item(11, "B").Value = 2012
$ws.Cells.Item(11, "C").Value = -0.2075757021743008
$ws.Cells.Item(11, "D").Value = 2013
$ws.Cells.Item(11, "E").Value = -0.2793004163246238
$ws.Cells.Item(12, "A").Value = 41409
$ws.Cells.Item(12, "B").Value = 2013
$ws.Cells.Item(12, "C").Value = 0.05915234751026066
$ws.Cells.Item(12, "D").Value = 2014
$ws.Cells.Item(12, "E").Value = 0.04624521867206965
$ws.Cells.Item(13, "A").Value = 41592
$ws.Cells.Item(13, "B").Value = 2013
$ws.Cells.Item(13, "C").Value = 0.124712275190686
$ws.Cells.Item(13, "D").Value = 2014
$ws.Cells.Item(13, "E").Value = -0.119752617912039
$ws.Cells.Item(14, "A").Value = 41774
$ws.Cells.Item(14, "B").Value = 2014
$ws.Cells.Item(14, "C").Value = -0.2979029954603529
$ws.Cells.Item(14, "D").Value = 2015
$ws.Cells.Item(14, "E").Value = -0.1124510725819206
$ws.Cells.Item(15, "A").Value = 41957
$ws.Cells.Item(15, "B").Value = 2014
$ws.Cells.Item(15, "C").Value = -0.255298189276465
$ws.Cells.Item(15, "D").Value = 2015
$ws.Cells.Item(15, "E").Value = -0.05946205208092747
$ws.Cells.Item(16, "A").Value = 42137
$ws.Cells.Item(16, "B").Value = 2015
$ws.Cells.Item(16, "C").Value = 0.07317408757452348
$ws.Cells.Item(16, "D").Value = 2016
$ws.Cells.Item(16, "E").Value = -0.002181547367274828
$ws.Cells.Item(17, "A").Value = 42321
$ws.Cells.Item(17, "B").Value = 2015
$ws.Cells.Item(17, "C").Value = 0.07418514192796266
$ws.Cells.Item(17, "D").Value = 2016
$ws.Cells.Item(17, "E").Value = -0.001680662521774678
$ws.Cells.Item(18, "A").Value = 42503
$ws.Cells.Item(18, "B").Value = 2016
$ws.Cells.Item(18, "C").Value = -0.06188089372189953
$ws.Cells.Item(18, "D").Value = 2017
$ws.Cells.Item(18, "E").Value = -0.07932008107318644
$ws.Cells.Item(19, "A").Value = 42689
$ws.Cells.Item(19, "B").Value = 2016
$ws.Cells.Item(19, "C").Value = -0.07611406013281474
$ws.Cells.Item(19, "D").Value = 2017
$ws.Cells.Item(19, "E").Value = -0.1247901924724348
$ws.Cells.Item(20, "A").Value = 42867
$ws.Cells.Item(20, "B").Value = 2017
$ws.Cells.Item(20, "C").Value = -0.2199961235931358
$ws.Cells.Item(20, "D").Value = 2018
$ws.Cells.Item(20, "E").Value = -0.1022879117640763
$ws.Cells.Item(21, "A").Value = 43053
$ws.Cells.Item(21, "B").Value = 2017
$ws.Cells.Item(21, "C").Value = -0.191300579729714
$ws.Cells.Item(21, "D").Value = 2018
$ws.Cells.Item(21, "E").Value = -0.05219951976568327
$ws.Cells.Item(22, "A").Value = 43145
$ws.Cells.Item(22, "B").Value = 2018
$ws.Cells.Item(22, "C").Value = 0.0441865668729946
$ws.Cells.Item(22, "D").Value = 2019
$ws.Cells.Item(22, "E").Value = -0.07120909843567613
$ws.Cells.Item(23, "A").Value = 43235
$ws.Cells.Item(23, "B").Value = 2018
$ws.Cells.Item(23, "C").Value = 0.0882025545300813
$ws.Cells.Item(23, "D").Value = 2019
$ws.Cells.Item(23, "E").Value = -0.05573300569792217
$ws.Cells.Item(24, "A").Value = 43326
$ws.Cells.Item(24, "B").Value = 2018
$ws.Cells.Item(24, "C").Value = 0.1415113532986956
$ws.Cells.Item(24, "D").Value = 2019
$ws.Cells.Item(24, "E").Value = 0.02926805735909976
$ws.Cells.Item(25, "A").Value = 43418
$ws.Cells.Item(25, "B").Value = 2018
$ws.Cells.Item(25, "C").Value = 0.0970330232288763
$ws.Cells.Item(25, "D").Value = 2019
$ws.Cells.Item(25, "E").Value = -0.1345737582127748
$ws.Cells.Item(26, "A").Value = 43510
$ws.Cells.Item(26, "B").Value = 2019
$ws.Cells.Item(26, "C").Value = -0.5756287392657988
$ws.Cells.Item(26, "D").Value = 2020
$ws.Cells.Item(26, "E").Value = -0.2394607875814136
$ws.Cells.Item(27, "A").Value = 43600
$ws.Cells.Item(27, "B").Value = 2019
$ws.Cells.Item(27, "C").Value = -0.7844010209450802
$ws.Cells.Item(27, "D").Value = 2020
$ws.Cells.Item(27, "E").Value = -0.3786583343736716
$ws.Cells.Item(28, "A").Value = 43691
$ws.Cells.Item(28, "B").Value = 2019
$ws.Cells.Item(28, "C").Value = -0.6919146680131605
$ws.Cells.Item(28, "D").Value = 2020
$ws.Cells.Item(28, "E").Value = -0.244860729922769
$ws.Cells.Item(29, "A").Value = 43783
$ws.Cells.Item(29, "B").Value = 2019
$ws.Cells.Item(29, "C").Value = -0.7407518902333265
$ws.Cells.Item(29, "D").Value = 2020
$ws.Cells.Item(29, "E").Value = -0.4363737508290888
$ws.Cells.Item(30, "A").Value = 43875
$ws.Cells.Item(30, "B").Value = 2020
$ws.Cells.Item(30, "C").Value = -0.7124953797697064
$ws.Cells.Item(30, "D").Value = 2021
$ws.Cells.Item(30, "E").Value = -0.4617192974095352
$ws.Cells.Item(31, "A").Value = 43966
$ws.Cells.Item(31, "B").Value = 2020
$ws.Cells.Item(31, "C").Value = -0.1808804304865297
$ws.Cells.Item(31, "D").Value = 2021
$ws.Cells.Item(31, "E").Value = -0.1077309791980285
$ws.Cells.Item(32, "A").Value = 44068
$ws.Cells.Item(32, "B").Value = 2020
$ws.Cells.Item(32, "C").Value = 0.3056679541520335
$ws.Cells.Item(32, "D").Value = 2021
$ws.Cells.Item(32, "E").Value = 0.3245880452514394
$ws.Cells.Item(33, "A").Value = 44159
$ws.Cells.Item(33, "B").Value = 2020
$ws.Cells.Item(33, "C").Value = 0.3056679541520335
$ws.Cells.Item(33, "D").Value = 2021
$ws.Cells.Item(33, "E").Value = -0.514812792200714
$ws.Cells.Item(34, "A").Value = 44251
$ws.Cells.Item(34, "B").Value = 2021
$ws.Cells.Item(34, "C").Value = -0.8680533514735522
$ws.Cells.Item(34, "D").Value = 2022
$ws.Cells.Item(34, "E").Value = -0.5995895195426981
$ws.Cells.Item(35, "A").Value = 44341
$ws.Cells.Item(35, "B").Value = 2021
$ws.Cells.Item(35, "C").Value = -0.8769761459347714
$ws.Cells.Item(35, "D").Value = 2022
$ws.Cells.Item(35, "E").Value = -0.5354669478056073
$ws.Cells.Item(36, "A").Value = 44432
$ws.Cells.Item(36, "B").Value = 2021
$ws.Cells.Item(36, "C").Value = -1.388491535160907
$ws.Cells.Item(36, "D").Value = 2022
$ws.Cells.Item(36, "E").Value = -2.541003699199929
$ws.Cells.Item(37, "A").Value = 44525
$ws.Cells.Item(37, "B").Value = 2021
$ws.Cells.Item(37, "C").Value = -1.388491535160907
$ws.Cells.Item(37, "D").Value = 2022
$ws.Cells.Item(37, "E").Value = -2.321721165370549
$ws.Cells.Item(38, "A").Value = 44617
$ws.Cells.Item(38, "B").Value = 2022
$ws.Cells.Item(38, "C").Value = -1.867377038014506
$ws.Cells.Item(38, "D").Value = 2023
$ws.Cells.Item(38, "E").Value = -0.8537083312609495
$ws.Cells.Item(39, "A").Value = 44706
$ws.Cells.Item(39, "B").Value = 2022
$ws.Cells.Item(39, "C").Value = -1.719168896439693
$ws.Cells.Item(39, "D").Value = 2023
$ws.Cells.Item(39, "E").Value = -0.5989817782328322
$ws.Cells.Item(40, "A").Value = 44798
$ws.Cells.Item(40, "B").Value = 2022
$ws.Cells.Item(40, "C").Value = -1.678482969789596
$ws.Cells.Item(40, "D").Value = 2023
$ws.Cells.Item(40, "E").Value = -0.6229862770763095
$ws.Cells.Item(41, "A").Value = 44890
$ws.Cells.Item(41, "B").Value = 2022
$ws.Cells.Item(41, "C").Value = -1.678482969789596
$ws.Cells.Item(41, "D").Value = 2023
$ws.Cells.Item(41, "E").Value = -1.107351089172237
$ws.Cells.Item(42, "A").Value = 44981
$ws.Cells.Item(42, "B").Value = 2023
$ws.Cells.Item(42, "C").Value = -0.9016470784766528
$ws.Cells.Item(42, "D").Value = 2024
$ws.Cells.Item(42, "E").Value = -1.556833564431637
$ws.Cells.Item(43, "A").Value = 45071
$ws.Cells.Item(43, "B").Value = 2023
$ws.Cells.Item(43, "C").Value = -0.6902657121583777
$ws.Cells.Item(43, "D").Value = 2024
$ws.Cells.Item(43, "E").Value = -1.134712300966823
$ws.Cells.Item(44, "A").Value = 45163
$ws.Cells.Item(44, "B").Value = 2023
$ws.Cells.Item(44, "C").Value = -0.5999457276250508
$ws.Cells.Item(44, "D").Value = 2024
$ws.Cells.Item(44, "E").Value = -0.7253995615808195
$ws.Cells.Item(45, "A").Value = 45254
$ws.Cells.Item(45, "B").Value = 2023
$ws.Cells.Item(45, "C").Value = -0.5999457276250508
$ws.Cells.Item(45, "D").Value = 2024
$ws.Cells.Item(45, "E").Value = -0.4628630633218611
$ws.Cells.Item(46, "A").Value = 45345
$ws.Cells.Item(46, "B").Value = 2024
$ws.Cells.Item(46, "C").Value = -0.101315145211045
$ws.Cells.Item(46, "D").Value = 2025
$ws.Cells.Item(46, "E").Value = -0.1859728711464226
$ws.Cells.Item(47, "A").Value = 45436
$ws.Cells.Item(47, "B").Value = 2024
$ws.Cells.Item(47, "C").Value = -0.07874066250703748
$ws.Cells.Item(47, "D").Value = 2025
$ws.Cells.Item(47, "E").Value = -0.1533081888441812
$ws.Cells.Item(48, "A").Value = 45534
$ws.Cells.Item(48, "B").Value = 2024
$ws.Cells.Item(48, "C").Value = -0.05499271238530445
$ws.Cells.Item(48, "D").Value = 2025
$ws.Cells.Item(48, "E").Value = -0.1365897193907339
$ws.Cells.Item(49, "A").Value = 45618
$ws.Cells.Item(49, "B").Value = 2024
$ws.Cells.Item(49, "C").Value = -0.05499271238530445
$ws.Cells.Item(49, "D").Value = 2025
$ws.Cells.Item(49, "E").Value = -0.01934819856548309
$ws.Cells.Item(50, "A").Value = 45713
$ws.Cells.Item(50, "B").Value = 2025
$ws.Cells.Item(50, "C").Value = 0.4236015715998187
$ws.Cells.Item(50, "D").Value = 2026
$ws.Cells.Item(50, "E").Value = 0.05974827491124213
$ws.Cells.Item(51, "A").Value = 45800
$ws.Cells.Item(51, "B").Value = 2025
$ws.Cells.Item(51, "C").Value = 0.4305325812036687
$ws.Cells.Item(51, "D").Value = 2026
$ws.Cells.Item(51, "E").Value = 0.2136583044595852
$ws.Cells.Item(52, "A").Value = 45891
$ws.Cells.Item(52, "B").Value = 2025
$ws.Cells.Item(52, "C").Value = 0.4335297397760618
$ws.Cells.Item(52, "D").Value = 2026
$ws.Cells.Item(52, "E").Value = 0.2794570629465865

Write-Host "Naive forecaster data bugfix applied"
